$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the product-code values in row 3 (Kho TC SPDK Nam Dinh)
$ws.Range("C3").Value = "HH009K"
$ws.Range("D3").Value = "HH0057"
$ws.Range("F3").Value = "HH050-021"

# Update the saved cursor/selection position
$ws.Range("H15").Select()
